# Update the users table on Sheet1 with the new AT/NU manufacturer and
# authorised-rep identifiers ("177H9/177H10" -> "277H12"), and move the
# active selection from A8 to D8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = "Manufacturer277H12_AT"
$ws.Range("A4").Value = "AuthorisedRep277H12_AT"
$ws.Range("A6").Value = "Manufacturer277H12_NU"
$ws.Range("A7").Value = "AuthorisedRep277H12_NU"

$ws.Range("D8").Select()
